# Update the bigram frequency table (rows 2-81) to the new
# re-ranked counts/pairs per the latest "Add files via upload".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "('climat', 'chang')"
$ws.Cells.Item(2, 3).Value = 169
$ws.Cells.Item(3, 2).Value = "('�', 't')"
$ws.Cells.Item(3, 3).Value = 61
$ws.Cells.Item(4, 2).Value = "('global', 'warm')"
$ws.Cells.Item(4, 3).Value = 35
$ws.Cells.Item(5, 2).Value = "('�', 'r')"
$ws.Cells.Item(5, 3).Value = 32
$ws.Cells.Item(6, 2).Value = "('sea', 'level')"
$ws.Cells.Item(6, 3).Value = 20
$ws.Cells.Item(7, 2).Value = "('don', '�')"
$ws.Cells.Item(7, 3).Value = 20
$ws.Cells.Item(8, 2).Value = "('polit', 'statement')"
$ws.Cells.Item(8, 3).Value = 16
$ws.Cells.Item(9, 2).Value = "('al', 'gore')"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(10, 2).Value = "('�', 'm')"
$ws.Cells.Item(10, 3).Value = 15
$ws.Cells.Item(11, 2).Value = "('year', 'ago')"
$ws.Cells.Item(11, 3).Value = 14
$ws.Cells.Item(12, 2).Value = "('3rd', 'world')"
$ws.Cells.Item(12, 3).Value = 14
$ws.Cells.Item(13, 2).Value = "('video', 'game')"
$ws.Cells.Item(13, 3).Value = 13
$ws.Cells.Item(14, 2).Value = "('level', 'rise')"
$ws.Cells.Item(14, 3).Value = 13
$ws.Cells.Item(15, 2).Value = "('reddit', 'kotakuinact')"
$ws.Cells.Item(15, 3).Value = 13
$ws.Cells.Item(16, 2).Value = "('kotakuinact', 'comment')"
$ws.Cells.Item(16, 3).Value = 13
$ws.Cells.Item(17, 2).Value = "('ice', 'cap')"
$ws.Cells.Item(17, 3).Value = 12
$ws.Cells.Item(18, 2).Value = "('last', 'year')"
$ws.Cells.Item(18, 3).Value = 12
$ws.Cells.Item(19, 2).Value = "('man', 'make')"
$ws.Cells.Item(19, 3).Value = 11
$ws.Cells.Item(20, 2).Value = "('ice', 'age')"
$ws.Cells.Item(20, 3).Value = 9
$ws.Cells.Item(21, 2).Value = "('hockey', 'stick')"
$ws.Cells.Item(21, 3).Value = 9
$ws.Cells.Item(22, 2).Value = "('late', 'game')"
$ws.Cells.Item(22, 3).Value = 9
$ws.Cells.Item(23, 2).Value = "('year', 'year')"
$ws.Cells.Item(23, 3).Value = 9
$ws.Cells.Item(24, 2).Value = "('gather', 'storm')"
$ws.Cells.Item(24, 3).Value = 8
$ws.Cells.Item(25, 2).Value = "('co2', 'emiss')"
$ws.Cells.Item(25, 3).Value = 8
$ws.Cells.Item(26, 2).Value = "('carbon', 'emiss')"
$ws.Cells.Item(26, 3).Value = 8
$ws.Cells.Item(27, 2).Value = "('nasa', 'gov')"
$ws.Cells.Item(27, 3).Value = 8
$ws.Cells.Item(28, 2).Value = "('pari', 'agreement')"
$ws.Cells.Item(28, 3).Value = 8
$ws.Cells.Item(29, 2).Value = "('specif', 'heat')"
$ws.Cells.Item(29, 3).Value = 8
$ws.Cells.Item(30, 2).Value = "('power', 'plant')"
$ws.Cells.Item(30, 3).Value = 7
$ws.Cells.Item(31, 2).Value = "('peopl', 'think')"
$ws.Cells.Item(31, 3).Value = 7
$ws.Cells.Item(32, 2).Value = "('co2', 'level')"
$ws.Cells.Item(32, 3).Value = 7
$ws.Cells.Item(33, 2).Value = "('black', 'peopl')"
$ws.Cells.Item(33, 3).Value = 7
$ws.Cells.Item(34, 2).Value = "('low', 'iq')"
$ws.Cells.Item(34, 3).Value = 7
$ws.Cells.Item(35, 2).Value = "('coal', 'oil')"
$ws.Cells.Item(35, 3).Value = 7
$ws.Cells.Item(36, 2).Value = "('�', '�')"
$ws.Cells.Item(36, 3).Value = 7
$ws.Cells.Item(37, 2).Value = "('�', 'll')"
$ws.Cells.Item(37, 3).Value = 7
$ws.Cells.Item(38, 2).Value = "('polit', 'issu')"
$ws.Cells.Item(38, 3).Value = 6
$ws.Cells.Item(39, 2).Value = "('melt', 'ice')"
$ws.Cells.Item(39, 3).Value = 6
$ws.Cells.Item(40, 2).Value = "('effect', 'climat')"
$ws.Cells.Item(40, 3).Value = 6
$ws.Cells.Item(41, 2).Value = "('carbon', 'dioxid')"
$ws.Cells.Item(41, 3).Value = 6
$ws.Cells.Item(42, 2).Value = "('realli', 'want')"
$ws.Cells.Item(42, 3).Value = 6
$ws.Cells.Item(43, 2).Value = "('lord', 'believ')"
$ws.Cells.Item(43, 3).Value = 6
$ws.Cells.Item(44, 2).Value = "('chang', 'polit')"
$ws.Cells.Item(44, 3).Value = 6
$ws.Cells.Item(45, 2).Value = "('degre', 'celsius')"
$ws.Cells.Item(45, 3).Value = 6
$ws.Cells.Item(46, 2).Value = "('climat', 'scientist')"
$ws.Cells.Item(46, 3).Value = 6
$ws.Cells.Item(47, 2).Value = "('climat', 'scienc')"
$ws.Cells.Item(47, 3).Value = 6
$ws.Cells.Item(48, 2).Value = "('global', 'climat')"
$ws.Cells.Item(48, 3).Value = 6
$ws.Cells.Item(49, 2).Value = "('greenhous', 'effect')"
$ws.Cells.Item(49, 3).Value = 6
$ws.Cells.Item(50, 2).Value = "('bell', 'curv')"
$ws.Cells.Item(50, 3).Value = 6
$ws.Cells.Item(51, 2).Value = "('tile', 'flood')"
$ws.Cells.Item(51, 3).Value = 6
$ws.Cells.Item(52, 2).Value = "('see', '�')"
$ws.Cells.Item(52, 3).Value = 6
$ws.Cells.Item(53, 2).Value = "('publish', 'report')"
$ws.Cells.Item(53, 3).Value = 5
$ws.Cells.Item(54, 2).Value = "('carbon', 'pollut')"
$ws.Cells.Item(54, 3).Value = 5
$ws.Cells.Item(55, 2).Value = "('global', 'catastroph')"
$ws.Cells.Item(55, 3).Value = 5
$ws.Cells.Item(56, 2).Value = "('thing', 'happen')"
$ws.Cells.Item(56, 3).Value = 5
$ws.Cells.Item(57, 2).Value = "('mind', 'worm')"
$ws.Cells.Item(57, 3).Value = 5
$ws.Cells.Item(58, 2).Value = "('climat', 'model')"
$ws.Cells.Item(58, 3).Value = 5
$ws.Cells.Item(59, 2).Value = "('natur', 'disast')"
$ws.Cells.Item(59, 3).Value = 5
$ws.Cells.Item(60, 2).Value = "('fossil', 'fuel')"
$ws.Cells.Item(60, 3).Value = 5
$ws.Cells.Item(61, 2).Value = "('chang', 'thing')"
$ws.Cells.Item(61, 3).Value = 5
$ws.Cells.Item(62, 2).Value = "('say', 'thing')"
$ws.Cells.Item(62, 3).Value = 5
$ws.Cells.Item(63, 2).Value = "('make', 'sen')"
$ws.Cells.Item(63, 3).Value = 5
$ws.Cells.Item(64, 2).Value = "('human', 'hive')"
$ws.Cells.Item(64, 3).Value = 5
$ws.Cells.Item(65, 2).Value = "('giss', 'nasa')"
$ws.Cells.Item(65, 3).Value = 5
$ws.Cells.Item(66, 2).Value = "('believ', 'climat')"
$ws.Cells.Item(66, 3).Value = 5
$ws.Cells.Item(67, 2).Value = "('go', 'back')"
$ws.Cells.Item(67, 3).Value = 5
$ws.Cells.Item(68, 2).Value = "('chang', 'real')"
$ws.Cells.Item(68, 3).Value = 5
$ws.Cells.Item(69, 2).Value = "('high', 'co2')"
$ws.Cells.Item(69, 3).Value = 5
$ws.Cells.Item(70, 2).Value = "('feel', 'good')"
$ws.Cells.Item(70, 3).Value = 5
$ws.Cells.Item(71, 2).Value = "('green', 'tech')"
$ws.Cells.Item(71, 3).Value = 5
$ws.Cells.Item(72, 2).Value = "('solar', 'panel')"
$ws.Cells.Item(72, 3).Value = 5
$ws.Cells.Item(73, 2).Value = "('long', 'time')"
$ws.Cells.Item(73, 3).Value = 5
$ws.Cells.Item(74, 2).Value = "('chang', 'climat')"
$ws.Cells.Item(74, 3).Value = 5
$ws.Cells.Item(75, 2).Value = "('polar', 'bear')"
$ws.Cells.Item(75, 3).Value = 5
$ws.Cells.Item(76, 2).Value = "('grow', 'wheat')"
$ws.Cells.Item(76, 3).Value = 5
$ws.Cells.Item(77, 2).Value = "('get', 'grip')"
$ws.Cells.Item(77, 3).Value = 5
$ws.Cells.Item(78, 2).Value = "('say', 'year')"
$ws.Cells.Item(78, 3).Value = 5
$ws.Cells.Item(79, 2).Value = "('take', 'much')"
$ws.Cells.Item(79, 3).Value = 5
$ws.Cells.Item(80, 2).Value = "('use', 'coal')"
$ws.Cells.Item(80, 3).Value = 5
$ws.Cells.Item(81, 2).Value = "('make', 'think')"
$ws.Cells.Item(81, 3).Value = 5
